{"js": "const doc = context.document;\nconst body = doc.body;\n\n// 1) Bold + italic + underline the word \"polimorfizam)\" in Zadatak 2.\nconst polyResults = body.search(\"polimorfizam)\", { matchCase: true });\npolyResults.load(\"text\");\nawait context.sync();\nconst polyRange = polyResults.items[0];\npolyRange.font.bold = true;\npolyRange.font.italic = true;\npolyRange.font.underline = \"Single\";\nawait context.sync();\n\n// 2) Insert \" - samo za RAM\" right before the closing parenthesis of the\n//    \"Napraviti interfejs...\" sentence (Zadatak 3). Restrict the search to\n//    that specific paragraph so we don't touch the similar word later on.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\nlet targetPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Napraviti interfejs\") !== -1) {\n    targetPara = paragraphs.items[i];\n    break;\n  }\n}\nconst pozResults = targetPara.search(\"pozajmice\", { matchCase: true });\npozResults.load(\"text\");\nawait context.sync();\npozResults.items[0].insertText(\" - samo za RAM\", \"End\");\nawait context.sync();\n\n// 3) Move the \"_GoBack\" bookmark (Word's marker for the most recent edit\n//    location) from its old spot at the very end of the document to the\n//    place we just edited in step 1 (between \"prodavnic\" and \"e.\").\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst gbResults = body.search(\"zarada prodavnic\", { matchCase: true });\ngbResults.load(\"text\");\nawait context.sync();\nconst gbRange = gbResults.items[0];\ngbRange.getRange(\"End\").insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the old \"_GoBack\" bookmark (Word will re-create it at the\n#    location of the most recent edit once we finish making changes below).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Bold + italic + underline the word \"polimorfizam)\" in Zadatak 2.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"polimorfizam)\")\nif ($found) {\n    $rng.Bold = 1\n    $rng.Italic = 1\n    $rng.Font.Underline = 1\n}\n\n# 3) Insert \" - samo za RAM\" right before the closing parenthesis in the\n#    \"Napraviti interfejs...\" sentence (Zadatak 3).\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$found2 = $rng2.Find.Execute(\"pozajmice)\")\nif ($found2) {\n    $insertPos = $rng2.Start + 9\n    $insertRange = $d.Range($insertPos, $insertPos)\n    $insertRange.InsertAfter(\" - samo za RAM\")\n}\n\n# 4) Re-create the \"_GoBack\" bookmark at the location of the edit made in\n#    step 2 (between \"prodavnic\" and \"e.\"), mirroring Word's behaviour of\n#    tracking the last editing location.\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$found3 = $rng3.Find.Execute(\"zarada prodavnice\")\nif ($found3) {\n    $bmPos = $rng3.Start + 16\n    $bmTarget = $d.Range($bmPos, $bmPos)\n    $d.Bookmarks.Add(\"_GoBack\", $bmTarget)\n}\n\nWrite-Output \"done\"\n"}
